{"js": "// Update the briefing document: refresh the send timestamp and fill in\n// several \"N\u00e3o informado\" (or placeholder) answers with real values.\n\n// 1) Header timestamp: \"Data de envio: 23/06/2025, 19:26:37\" -> \"...19:37:59\"\nconst tsResults = context.document.body.search(\n  \"Data de envio: 23/06/2025, 19:26:37\",\n  { matchCase: true }\n);\ntsResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < tsResults.items.length; i++) {\n  tsResults.items[i].insertText(\n    \"Data de envio: 23/06/2025, 19:37:59\",\n    \"Replace\"\n  );\n}\n\n// 2) Label -> [oldValue, newValue] pairs. Each label/value lives in the\n// same paragraph, separated by a manual line break, e.g.:\n//   \"Nome completo da empresa:<br> N\u00e3o informado\"\nconst fieldChanges = [\n  [\"Nome completo da empresa:\", \" N\u00e3o informado\", \" Teste\"],\n  [\"Lista de produtos/servi\u00e7os:\", \" N\u00e3o informado\", \" Ser\u00e1 que funciona?\"],\n  [\"Canais de compra:\", \" N\u00e3o informado\", \" Testando \"],\n  [\"Como clientes procuram:\", \" N\u00e3o informado\", \" Teste\"],\n  [\"Concorrentes diretos:\", \" N\u00e3o informado\", \" Testando \"],\n  [\"Identidade visual:\", \" N\u00e3o informado\", \" N\u00e3o possui logo\"],\n  [\"Redes sociais:\", \" 5654\", \" N\u00e3o informado\"],\n  [\"Site pr\u00f3prio:\", \" N\u00e3o informado\", \" Hahshss\"],\n  [\"Atributos do neg\u00f3cio:\", \" N\u00e3o informado\", \" Meme\"],\n  [\"Avalia\u00e7\u00f5es online recebidas:\", \" N\u00e3o informado\", \" Sim\"],\n  [\"Onde recebeu avalia\u00e7\u00f5es:\", \" dasdsa\", \" N\u00e3o informado\"],\n  [\"Objetivo principal:\", \" N\u00e3o informado\", \" Melhorar reputa\u00e7\u00e3o online\"],\n  [\"Google Ads futuro:\", \" Talvez\", \" N\u00e3o informado\"],\n  [\"Respons\u00e1vel pela gest\u00e3o:\", \" dsadasdasdsa\", \" Aaaa\"],\n  [\"Tentativa anterior GMB:\", \" 4654\", \" N\u00e3o informado\"],\n  [\"Informa\u00e7\u00f5es a ocultar:\", \" dasdsadsa\", \" N\u00e3o informado\"],\n  [\"Detalhes importantes:\", \" dasdsadsa\", \" N\u00e3o informado\"],\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const [label, oldValue, newValue] of fieldChanges) {\n  let target = null;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.indexOf(label) === 0) {\n      target = paragraphs.items[i];\n      break;\n    }\n  }\n  if (!target) {\n    continue;\n  }\n  const valueResults = target.getRange().search(oldValue, { matchCase: true });\n  valueResults.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < valueResults.items.length; i++) {\n    valueResults.items[i].insertText(newValue, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Header timestamp: \"Data de envio: 23/06/2025, 19:26:37\" -> \"...19:37:59\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Data de envio: 23/06/2025, 19:26:37\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Data de envio: 23/06/2025, 19:37:59\"\n$find.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2) | Out-Null\n\n# 2) Label -> old/new value pairs. Each label/value lives in the same\n# paragraph, separated by a manual line break, e.g.:\n#   \"Nome completo da empresa:<br> N\u00e3o informado\"\n$fieldChanges = @(\n    @(\"Nome completo da empresa:\", \" N\u00e3o informado\", \" Teste\"),\n    @(\"Lista de produtos/servi\u00e7os:\", \" N\u00e3o informado\", \" Ser\u00e1 que funciona?\"),\n    @(\"Canais de compra:\", \" N\u00e3o informado\", \" Testando \"),\n    @(\"Como clientes procuram:\", \" N\u00e3o informado\", \" Teste\"),\n    @(\"Concorrentes diretos:\", \" N\u00e3o informado\", \" Testando \"),\n    @(\"Identidade visual:\", \" N\u00e3o informado\", \" N\u00e3o possui logo\"),\n    @(\"Redes sociais:\", \" 5654\", \" N\u00e3o informado\"),\n    @(\"Site pr\u00f3prio:\", \" N\u00e3o informado\", \" Hahshss\"),\n    @(\"Atributos do neg\u00f3cio:\", \" N\u00e3o informado\", \" Meme\"),\n    @(\"Avalia\u00e7\u00f5es online recebidas:\", \" N\u00e3o informado\", \" Sim\"),\n    @(\"Onde recebeu avalia\u00e7\u00f5es:\", \" dasdsa\", \" N\u00e3o informado\"),\n    @(\"Objetivo principal:\", \" N\u00e3o informado\", \" Melhorar reputa\u00e7\u00e3o online\"),\n    @(\"Google Ads futuro:\", \" Talvez\", \" N\u00e3o informado\"),\n    @(\"Respons\u00e1vel pela gest\u00e3o:\", \" dsadasdasdsa\", \" Aaaa\"),\n    @(\"Tentativa anterior GMB:\", \" 4654\", \" N\u00e3o informado\"),\n    @(\"Informa\u00e7\u00f5es a ocultar:\", \" dasdsadsa\", \" N\u00e3o informado\"),\n    @(\"Detalhes importantes:\", \" dasdsadsa\", \" N\u00e3o informado\")\n)\n\nforeach ($change in $fieldChanges) {\n    $label = $change[0]\n    $oldValue = $change[1]\n    $newValue = $change[2]\n    foreach ($p in $d.Paragraphs) {\n        $t = $p.Range.Text\n        if ($t.StartsWith($label)) {\n            $rng = $p.Range\n            $pfind = $rng.Find\n            $pfind.ClearFormatting()\n            $pfind.Text = $oldValue\n            $pfind.Replacement.ClearFormatting()\n            $pfind.Replacement.Text = $newValue\n            $pfind.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2) | Out-Null\n            break\n        }\n    }\n}\n"}
